# Update the "Förändrad" (Changed) date column (C) for rows 2-41
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
